$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

    $ws.Cells.Item(2, 7).Value = 18.629453
    $ws.Cells.Item(2, 8).Value = 55.888359
    $ws.Cells.Item(2, 9).Value = 0.07116572597273459
    $ws.Cells.Item(2, 10).Value = 0.07977938953593823
    $ws.Cells.Item(2, 13).Value = 18.10884
    $ws.Cells.Item(2, 14).Value = 54.32652
    $ws.Cells.Item(2, 15).Value = 0.02835750805894763
    $ws.Cells.Item(2, 16).Value = 0.02862140385105863
    $ws.Cells.Item(2, 17).Value = 337.35778366452
    $ws.Cells.Item(2, 18).Value = 3036.22005298068
    $ws.Cells.Item(2, 19).Value = 0.00201808264779268
    $ws.Cells.Item(2, 20).Value = 0.002283398126899009
    $ws.Cells.Item(3, 7).Value = 18.629453
    $ws.Cells.Item(3, 8).Value = 55.888359
    $ws.Cells.Item(3, 9).Value = 0.07116572597273459
    $ws.Cells.Item(3, 10).Value = 0.07977938953593823
    $ws.Cells.Item(3, 15).Value = 0.0007420423715060531
    $ws.Cells.Item(3, 16).Value = 0.0007489478393278889
    $ws.Cells.Item(3, 17).Value = 8.827777437850667
    $ws.Cells.Item(3, 18).Value = 79.44999694065599
    $ws.Cells.Item(3, 19).Value = 0.00005280798407075789
    $ws.Cells.Item(3, 20).Value = 0.00005975060141583893
    $ws.Cells.Item(4, 7).Value = 18.629453
    $ws.Cells.Item(4, 8).Value = 55.888359
    $ws.Cells.Item(4, 9).Value = 0.07116572597273459
    $ws.Cells.Item(4, 10).Value = 0.07977938953593823
    $ws.Cells.Item(4, 13).Value = 305.722738
    $ws.Cells.Item(4, 14).Value = 917.168214
    $ws.Cells.Item(4, 15).Value = 0.4787460161246407
    $ws.Cells.Item(4, 16).Value = 0.4832012404300546
    $ws.Cells.Item(4, 17).Value = 5695.447378602315
    $ws.Cells.Item(4, 18).Value = 51259.02640742083
    $ws.Cells.Item(4, 19).Value = 0.03407030779406455
    $ws.Cells.Item(4, 20).Value = 0.03854949998451787
    $ws.Cells.Item(5, 7).Value = 18.629453
    $ws.Cells.Item(5, 8).Value = 55.888359
    $ws.Cells.Item(5, 9).Value = 0.07116572597273459
    $ws.Cells.Item(5, 10).Value = 0.07977938953593823
    $ws.Cells.Item(5, 13).Value = 17.6638495
    $ws.Cells.Item(5, 14).Value = 35.327699
    $ws.Cells.Item(5, 15).Value = 0.02766067592116823
    $ws.Cells.Item(5, 16).Value = 0.01861205798213543
    $ws.Cells.Item(5, 17).Value = 329.0678540593235
    $ws.Cells.Item(5, 18).Value = 1974.407124355941
    $ws.Cells.Item(5, 19).Value = 0.001968492082826476
    $ws.Cells.Item(5, 20).Value = 0.001484858623822251
    $ws.Cells.Item(6, 7).Value = 18.629453
    $ws.Cells.Item(6, 8).Value = 55.888359
    $ws.Cells.Item(6, 9).Value = 0.07116572597273459
    $ws.Cells.Item(6, 10).Value = 0.07977938953593823
    $ws.Cells.Item(6, 13).Value = 296.6213786666667
    $ws.Cells.Item(6, 14).Value = 889.864136
    $ws.Cells.Item(6, 15).Value = 0.4644937575237376
    $ws.Cells.Item(6, 16).Value = 0.4688163498974233
    $ws.Cells.Item(6, 17).Value = 5525.894032665869
    $ws.Cells.Item(6, 18).Value = 49733.04629399283
    $ws.Cells.Item(6, 19).Value = 0.03305603546398013
    $ws.Cells.Item(6, 20).Value = 0.03740188219928325
    $ws.Cells.Item(7, 9).Value = 0.2779443552245922
    $ws.Cells.Item(7, 10).Value = 0.3115858186182692
    $ws.Cells.Item(7, 13).Value = 18.10884
    $ws.Cells.Item(7, 14).Value = 54.32652
    $ws.Cells.Item(7, 15).Value = 0.02835750805894763
    $ws.Cells.Item(7, 16).Value = 0.02862140385105863
    $ws.Cells.Item(7, 17).Value = 1317.58217005412
    $ws.Cells.Item(7, 18).Value = 11858.23953048708
    $ws.Cells.Item(7, 19).Value = 0.007881809293220377
    $ws.Cells.Item(7, 20).Value = 0.008918023548936187
    $ws.Cells.Item(8, 9).Value = 0.2779443552245922
    $ws.Cells.Item(8, 10).Value = 0.3115858186182692
    $ws.Cells.Item(8, 15).Value = 0.0007420423715060531
    $ws.Cells.Item(8, 16).Value = 0.0007489478393278889
    $ws.Cells.Item(8, 19).Value = 0.0002062464884975772
    $ws.Cells.Item(8, 20).Value = 0.0002333615256193642
    $ws.Cells.Item(9, 9).Value = 0.2779443552245922
    $ws.Cells.Item(9, 10).Value = 0.3115858186182692
    $ws.Cells.Item(9, 13).Value = 305.722738
    $ws.Cells.Item(9, 14).Value = 917.168214
    $ws.Cells.Item(9, 15).Value = 0.4787460161246407
    $ws.Cells.Item(9, 16).Value = 0.4832012404300546
    $ws.Cells.Item(9, 17).Value = 22244.0989355987
    $ws.Cells.Item(9, 18).Value = 200196.8904203883
    $ws.Cells.Item(9, 19).Value = 0.1330647527681054
    $ws.Cells.Item(9, 20).Value = 0.1505586540567617
    $ws.Cells.Item(10, 9).Value = 0.2779443552245922
    $ws.Cells.Item(10, 10).Value = 0.3115858186182692
    $ws.Cells.Item(10, 13).Value = 17.6638495
    $ws.Cells.Item(10, 14).Value = 35.327699
    $ws.Cells.Item(10, 15).Value = 0.02766067592116823
    $ws.Cells.Item(10, 16).Value = 0.01861205798213543
    $ws.Cells.Item(10, 17).Value = 1285.20507971352
    $ws.Cells.Item(10, 18).Value = 7711.23047828112
    $ws.Cells.Item(10, 19).Value = 0.007688128733985507
    $ws.Cells.Item(10, 20).Value = 0.005799253322534361
    $ws.Cells.Item(11, 9).Value = 0.2779443552245922
    $ws.Cells.Item(11, 10).Value = 0.3115858186182692
    $ws.Cells.Item(11, 13).Value = 296.6213786666667
    $ws.Cells.Item(11, 14).Value = 889.864136
    $ws.Cells.Item(11, 15).Value = 0.4644937575237376
    $ws.Cells.Item(11, 16).Value = 0.4688163498974233
    $ws.Cells.Item(11, 17).Value = 21581.89258881693
    $ws.Cells.Item(11, 18).Value = 194237.0332993523
    $ws.Cells.Item(11, 19).Value = 0.1291034179407833
    $ws.Cells.Item(11, 20).Value = 0.1460765261644176
    $ws.Cells.Item(12, 7).Value = 36.272704
    $ws.Cells.Item(12, 8).Value = 108.818112
    $ws.Cells.Item(12, 9).Value = 0.1385640959589159
    $ws.Cells.Item(12, 10).Value = 0.1553354348051864
    $ws.Cells.Item(12, 13).Value = 18.10884
    $ws.Cells.Item(12, 14).Value = 54.32652
    $ws.Cells.Item(12, 15).Value = 0.02835750805894763
    $ws.Cells.Item(12, 16).Value = 0.02862140385105863
    $ws.Cells.Item(12, 17).Value = 656.8565931033598
    $ws.Cells.Item(12, 18).Value = 5911.70933793024
    $ws.Cells.Item(12, 19).Value = 0.003929332467835751
    $ws.Cells.Item(12, 20).Value = 0.004445918211939031
    $ws.Cells.Item(13, 7).Value = 36.272704
    $ws.Cells.Item(13, 8).Value = 108.818112
    $ws.Cells.Item(13, 9).Value = 0.1385640959589159
    $ws.Cells.Item(13, 10).Value = 0.1553354348051864
    $ws.Cells.Item(13, 15).Value = 0.0007420423715060531
    $ws.Cells.Item(13, 16).Value = 0.0007489478393278889
    $ws.Cells.Item(13, 17).Value = 17.18823188104533
    $ws.Cells.Item(13, 18).Value = 154.694086929408
    $ws.Cells.Item(13, 19).Value = 0.0001028204303709462
    $ws.Cells.Item(13, 20).Value = 0.0001163381382684025
    $ws.Cells.Item(14, 7).Value = 36.272704
    $ws.Cells.Item(14, 8).Value = 108.818112
    $ws.Cells.Item(14, 9).Value = 0.1385640959589159
    $ws.Cells.Item(14, 10).Value = 0.1553354348051864
    $ws.Cells.Item(14, 13).Value = 305.722738
    $ws.Cells.Item(14, 14).Value = 917.168214
    $ws.Cells.Item(14, 15).Value = 0.4787460161246407
    $ws.Cells.Item(14, 16).Value = 0.4832012404300546
    $ws.Cells.Item(14, 17).Value = 11089.39038154355
    $ws.Cells.Item(14, 18).Value = 99804.51343389197
    $ws.Cells.Item(14, 19).Value = 0.0663370089182434
    $ws.Cells.Item(14, 20).Value = 0.07505827478060798
    $ws.Cells.Item(15, 7).Value = 36.272704
    $ws.Cells.Item(15, 8).Value = 108.818112
    $ws.Cells.Item(15, 9).Value = 0.1385640959589159
    $ws.Cells.Item(15, 10).Value = 0.1553354348051864
    $ws.Cells.Item(15, 13).Value = 17.6638495
    $ws.Cells.Item(15, 14).Value = 35.327699
    $ws.Cells.Item(15, 15).Value = 0.02766067592116823
    $ws.Cells.Item(15, 16).Value = 0.01861205798213543
    $ws.Cells.Item(15, 17).Value = 640.7155844140478
    $ws.Cells.Item(15, 18).Value = 3844.293506484287
    $ws.Cells.Item(15, 19).Value = 0.003832776552629229
    $ws.Cells.Item(15, 20).Value = 0.002891112119274349
    $ws.Cells.Item(16, 7).Value = 36.272704
    $ws.Cells.Item(16, 8).Value = 108.818112
    $ws.Cells.Item(16, 9).Value = 0.1385640959589159
    $ws.Cells.Item(16, 10).Value = 0.1553354348051864
    $ws.Cells.Item(16, 13).Value = 296.6213786666667
    $ws.Cells.Item(16, 14).Value = 889.864136
    $ws.Cells.Item(16, 15).Value = 0.4644937575237376
    $ws.Cells.Item(16, 16).Value = 0.4688163498974233
    $ws.Cells.Item(16, 17).Value = 10759.25946844791
    $ws.Cells.Item(16, 18).Value = 96833.33521603124
    $ws.Cells.Item(16, 19).Value = 0.06436215758983659
    $ws.Cells.Item(16, 20).Value = 0.07282379155509668
    $ws.Cells.Item(17, 7).Value = 84.7905925
    $ws.Cells.Item(17, 8).Value = 169.581185
    $ws.Cells.Item(17, 9).Value = 0.3239055956672912
    $ws.Cells.Item(17, 10).Value = 0.2420733701642771
    $ws.Cells.Item(17, 13).Value = 18.10884
    $ws.Cells.Item(17, 14).Value = 54.32652
    $ws.Cells.Item(17, 15).Value = 0.02835750805894763
    $ws.Cells.Item(17, 16).Value = 0.02862140385105863
    $ws.Cells.Item(17, 17).Value = 1535.4592730877
    $ws.Cells.Item(17, 18).Value = 9212.7556385262
    $ws.Cells.Item(17, 19).Value = 0.009185155539473444
    $ws.Cells.Item(17, 20).Value = 0.006928479689058582
    $ws.Cells.Item(18, 7).Value = 84.7905925
    $ws.Cells.Item(18, 8).Value = 169.581185
    $ws.Cells.Item(18, 9).Value = 0.3239055956672912
    $ws.Cells.Item(18, 10).Value = 0.2420733701642771
    $ws.Cells.Item(18, 15).Value = 0.0007420423715060531
    $ws.Cells.Item(18, 16).Value = 0.0007489478393278889
    $ws.Cells.Item(18, 17).Value = 40.17898321617333
    $ws.Cells.Item(18, 18).Value = 241.07389929704
    $ws.Cells.Item(18, 19).Value = 0.0002403516763530375
    $ws.Cells.Item(18, 20).Value = 0.0001813003275433556
    $ws.Cells.Item(19, 7).Value = 84.7905925
    $ws.Cells.Item(19, 8).Value = 169.581185
    $ws.Cells.Item(19, 9).Value = 0.3239055956672912
    $ws.Cells.Item(19, 10).Value = 0.2420733701642771
    $ws.Cells.Item(19, 13).Value = 305.722738
    $ws.Cells.Item(19, 14).Value = 917.168214
    $ws.Cells.Item(19, 15).Value = 0.4787460161246407
    $ws.Cells.Item(19, 16).Value = 0.4832012404300546
    $ws.Cells.Item(19, 17).Value = 25922.41209574226
    $ws.Cells.Item(19, 18).Value = 155534.4725744536
    $ws.Cells.Item(19, 19).Value = 0.1550685135261943
    $ws.Cells.Item(19, 20).Value = 0.1169701527384625
    $ws.Cells.Item(20, 7).Value = 84.7905925
    $ws.Cells.Item(20, 8).Value = 169.581185
    $ws.Cells.Item(20, 9).Value = 0.3239055956672912
    $ws.Cells.Item(20, 10).Value = 0.2420733701642771
    $ws.Cells.Item(20, 13).Value = 17.6638495
    $ws.Cells.Item(20, 14).Value = 35.327699
    $ws.Cells.Item(20, 15).Value = 0.02766067592116823
    $ws.Cells.Item(20, 16).Value = 0.01861205798213543
    $ws.Cells.Item(20, 17).Value = 1497.728264935829
    $ws.Cells.Item(20, 18).Value = 5990.913059743314
    $ws.Cells.Item(20, 19).Value = 0.008959447710805894
    $ws.Cells.Item(20, 20).Value = 0.004505483601428458
    $ws.Cells.Item(21, 7).Value = 84.7905925
    $ws.Cells.Item(21, 8).Value = 169.581185
    $ws.Cells.Item(21, 9).Value = 0.3239055956672912
    $ws.Cells.Item(21, 10).Value = 0.2420733701642771
    $ws.Cells.Item(21, 13).Value = 296.6213786666667
    $ws.Cells.Item(21, 14).Value = 889.864136
    $ws.Cells.Item(21, 15).Value = 0.4644937575237376
    $ws.Cells.Item(21, 16).Value = 0.4688163498974233
    $ws.Cells.Item(21, 17).Value = 25150.70244531353
    $ws.Cells.Item(21, 18).Value = 150904.2146718812
    $ws.Cells.Item(21, 19).Value = 0.1504521272144645
    $ws.Cells.Item(21, 20).Value = 0.1134879538077842
    $ws.Cells.Item(22, 7).Value = 49.32382433333333
    $ws.Cells.Item(22, 8).Value = 147.971473
    $ws.Cells.Item(22, 9).Value = 0.1884202271764661
    $ws.Cells.Item(22, 10).Value = 0.2112259868763291
    $ws.Cells.Item(22, 13).Value = 18.10884
    $ws.Cells.Item(22, 14).Value = 54.32652
    $ws.Cells.Item(22, 15).Value = 0.02835750805894763
    $ws.Cells.Item(22, 16).Value = 0.02862140385105863
    $ws.Cells.Item(22, 17).Value = 893.1972430404398
    $ws.Cells.Item(22, 18).Value = 8038.775187363959
    $ws.Cells.Item(22, 19).Value = 0.005343128110625381
    $ws.Cells.Item(22, 20).Value = 0.006045584274225825
    $ws.Cells.Item(23, 7).Value = 49.32382433333333
    $ws.Cells.Item(23, 8).Value = 147.971473
    $ws.Cells.Item(23, 9).Value = 0.1884202271764661
    $ws.Cells.Item(23, 10).Value = 0.2112259868763291
    $ws.Cells.Item(23, 15).Value = 0.0007420423715060531
    $ws.Cells.Item(23, 16).Value = 0.0007489478393278889
    $ws.Cells.Item(23, 17).Value = 23.37265316369244
    $ws.Cells.Item(23, 18).Value = 210.353878473232
    $ws.Cells.Item(23, 19).Value = 0.0001398157922137342
    $ws.Cells.Item(23, 20).Value = 0.0001581972464809277
    $ws.Cells.Item(24, 7).Value = 49.32382433333333
    $ws.Cells.Item(24, 8).Value = 147.971473
    $ws.Cells.Item(24, 9).Value = 0.1884202271764661
    $ws.Cells.Item(24, 10).Value = 0.2112259868763291
    $ws.Cells.Item(24, 13).Value = 305.722738
    $ws.Cells.Item(24, 14).Value = 917.168214
    $ws.Cells.Item(24, 15).Value = 0.4787460161246407
    $ws.Cells.Item(24, 16).Value = 0.4832012404300546
    $ws.Cells.Item(24, 17).Value = 15079.41462381769
    $ws.Cells.Item(24, 18).Value = 135714.7316143592
    $ws.Cells.Item(24, 19).Value = 0.09020543311803289
    $ws.Cells.Item(24, 20).Value = 0.1020646588697046
    $ws.Cells.Item(25, 7).Value = 49.32382433333333
    $ws.Cells.Item(25, 8).Value = 147.971473
    $ws.Cells.Item(25, 9).Value = 0.1884202271764661
    $ws.Cells.Item(25, 10).Value = 0.2112259868763291
    $ws.Cells.Item(25, 13).Value = 17.6638495
    $ws.Cells.Item(25, 14).Value = 35.327699
    $ws.Cells.Item(25, 15).Value = 0.02766067592116823
    $ws.Cells.Item(25, 16).Value = 0.01861205798213543
    $ws.Cells.Item(25, 17).Value = 871.2486097884378
    $ws.Cells.Item(25, 18).Value = 5227.491658730626
    $ws.Cells.Item(25, 19).Value = 0.005211830840921124
    $ws.Cells.Item(25, 20).Value = 0.003931350315076015
    $ws.Cells.Item(26, 7).Value = 49.32382433333333
    $ws.Cells.Item(26, 8).Value = 147.971473
    $ws.Cells.Item(26, 9).Value = 0.1884202271764661
    $ws.Cells.Item(26, 10).Value = 0.2112259868763291
    $ws.Cells.Item(26, 13).Value = 296.6213786666667
    $ws.Cells.Item(26, 14).Value = 889.864136
    $ws.Cells.Item(26, 15).Value = 0.4644937575237376
    $ws.Cells.Item(26, 16).Value = 0.4688163498974233
